$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 34 appended at the bottom of the results table.
# Column A holds a numeric-looking value ("202474") that must stay text
# (matching the rest of column A), so mark it as Text before entering it.
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "202474"
$ws.Range("B34").Value = "100/100"
$ws.Range("C34").Value = "03:14"
